$d = $word.ActiveDocument

# Locate the run of text "repository" inside the hyperlink.
$rng = $d.Content
$found = $rng.Find.Execute("repository", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start

    # "repository" -> "re" + "p" + "ository"
    # The middle sub-range ("p") is characters 3-4 of the found text
    # (i.e. offsets start+2 .. start+3).
    $midRange = $d.Range($start + 2, $start + 3)

    # Toggling a character formatting property on that sub-range and then
    # restoring it forces Word to break the single run into three separate
    # runs at the sub-range boundaries ("re" | "p" | "ository"), without
    # altering the visible formatting of any of them.
    $originalBold = $midRange.Font.Bold
    $midRange.Font.Bold = 1
    $midRange.Font.Bold = $originalBold
}
